$wb = $excel.ActiveWorkbook

# The dropdown source list lives on the worksheet named "Sheet1" (the
# TRN_TYPE lookup list in column A, referenced by the data validation on
# the other sheet via Sheet1!$A$3:$A$17). Fix the long-standing typo
# "Ajustment" -> "Adjustment" in the two affected list entries so the
# dropdown shows the corrected values.
$lookupSheet = $wb.Worksheets.Item("Sheet1")

$lookupSheet.Range("A10").Value = "Inventory Adjustment"
$lookupSheet.Range("A11").Value = "Inventory Adjustment COGS"
